$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current "Total" row (row 12); the Total row
# (and its SUM formulas) shift down to row 13.
$ws.Rows.Item(12).Insert()

# Populate the new row 12 with another working-hours entry, following the
# same pattern as the existing data rows (set values/formulas first).
$ws.Range("A12").Value = 45274
$ws.Range("B12").Value = 0.416666666666667
$ws.Range("C12").Value = 0.5
$ws.Range("D12").Formula = "=(C12<B12)+C12-B12"
$ws.Range("E12").Value = 10
$ws.Range("F12").Formula = "=(D12*24)*E12"

# Copy the formatting (number formats / borders / fonts) from the row above
# onto the new row so it matches the rest of the data rows. Do this after
# writing the values/formulas so the pasted number formats stick (entering a
# formula otherwise re-infers a number format from its precedents).
$ws.Range("A11:F11").Copy()
$ws.Range("A12:F12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Extend the Total row (now shifted down to row 13) so its sums include the
# newly inserted row.
$ws.Range("D13").Formula = "=SUM(D2:D12)"
$ws.Range("F13").Formula = "=SUM(F2:F12)"

# Leave the selection where the authored workbook left it.
$ws.Range("A13").Select()
